$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: age 32 -> "33" stored as text, and clear the explicit
#     (Hyperlink-ish) styling that B2/C2 previously carried ---
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "33"
$ws.Range("F2").Style = "Normal"

$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Style = "Normal"

# --- Row 9: fill in the new "testing" staff record (admin, stage 1) ---
$ws.Range("A9").Style = "Normal"
$ws.Range("A9").Value = "244e4090-f72c-4ef3-8817-0f88efcf78ed"

$ws.Range("B9").Style = "Normal"
$ws.Range("B9").Value = "testing"

$ws.Range("C9").Value = "tester"
$ws.Range("D9").Value = "S"
$ws.Range("E9").Value = "M"

$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "43"
$ws.Range("F9").Style = "Normal"

$ws.Range("G9").Value = "NTU"

# --- Update the current selection to match the post-edit UI state ---
$ws.Range("D9").Select()
